# HourProcessor.cs: fields changed to dictionary. Compute Hours now working.
# Populate the newly-computed overtime-multiplier table in column C (rows 17-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# day/night * normal/special/regular overtime multipliers produced by the
# (now working) HourProcessor.ComputeHours().
$values = @(
    1,                      # C17
    1.3,                    # C18
    2,                      # C19
    1.1000000000000001,     # C20
    1.43,                   # C21
    2.2000000000000002,     # C22
    1.25,                   # C23
    1.69,                   # C24
    2.6,                    # C25
    1.375,                  # C26
    1.859,                  # C27
    2.86,                   # C28
    1.3,                    # C29
    1.5,                    # C30
    2.6,                    # C31
    1.43,                   # C32
    1.65,                   # C33
    2.86,                   # C34
    1.625,                  # C35
    1.95,                   # C36
    3.38,                   # C37
    1.859,                  # C38
    2.145,                  # C39
    3.718                   # C40
)

$row = 17
foreach ($v in $values) {
    $ws.Cells.Item($row, 3).Value = $v
    $row++
}

# Match the author's view state: zoomed to 95%, scrolled so row 2 is visible
# at the top, with the new table (C17:C40) selected.
$win.Zoom = 95
$ws.Range("C17:C40").Select()
